$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns I and J
# Copy formatting from the existing header cell (H1) so the new
# header cells share the same bold/bordered/centered style.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for columns I (I0) and J (IF), rows 2-31
$data = @(
    @(9, 9),
    @(9, 9),
    @(6, 6),
    @(9, 9),
    @(4, 5),
    @(7, 7),
    @(3, 4),
    @(9, 9),
    @(4, 4),
    @(1, 1),
    @(7, 7),
    @(3, 4),
    @(8, 9),
    @(7, 7),
    @(9, 9),
    @(6, 7),
    @(3, 3),
    @(6, 6),
    @(11, 11),
    @(7, 7),
    @(8, 8),
    @(5, 5),
    @(8, 8),
    @(8, 8),
    @(9, 9),
    @(8, 9),
    @(8, 9),
    @(8, 9),
    @(6, 6),
    @(8, 9)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
